$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores prices as literal text, e.g. "57.879.37" (dot-grouped
# thousands) or plain "526.01". Some of the new prices (e.g. "527.41", "0.999") are
# themselves valid numbers, so Excel would silently convert them to numeric cells on a
# plain `.Value =` assignment. Mark those specific cells as Text first so the literal
# string is preserved, then restore the "Normal" cell style so no visible/structural
# formatting change is left behind - only the cells stored value stays textual.
$textForceCells = @("D5", "D6", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.923.72'
$ws.Range("E2").Value = '  +2.56%  '
$ws.Range("D3").Value = '3.060.78'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '527.41'
$ws.Range("E5").Value = '  +6.35%  '
$ws.Range("D6").Value = '143.22'
$ws.Range("E6").Value = '  +6.41%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +5.32%  '
$ws.Range("E9").Value = '  +6.76%  '
$ws.Range("E10").Value = '  +8.47%  '
$ws.Range("E11").Value = '  +6.26%  '
$ws.Range("E12").Value = '  +2.51%  '
$ws.Range("D13").Value = '3.580.24'
$ws.Range("E13").Value = '  +2.58%  '
$ws.Range("D14").Value = '27.21'
$ws.Range("E14").Value = '  +8.78%  '
$ws.Range("D15").Value = '0.0000170'
$ws.Range("E15").Value = '  +16.89%  '
$ws.Range("D16").Value = '57.869.59'
$ws.Range("E16").Value = '  +2.52%  '
$ws.Range("E17").Value = '  +8.24%  '
$ws.Range("D18").Value = '3.057.69'
$ws.Range("E18").Value = '  +2.48%  '
$ws.Range("D19").Value = '13.08'
$ws.Range("E19").Value = '  +5.93%  '
$ws.Range("D20").Value = '8.12'
$ws.Range("E20").Value = '  +4.73%  '
$ws.Range("D21").Value = '339.35'
$ws.Range("E21").Value = '  +4.60%  '
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '0.506'
$ws.Range("E23").Value = '  +8.18%  '
$ws.Range("D24").Value = '65.09'
$ws.Range("E24").Value = '  +5.69%  '
$ws.Range("E25").Value = '  +6.15%  '
$ws.Range("D26").Value = '0.0₃0983'
$ws.Range("E26").Value = '  +9.98%  '
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").Value = '6.97'
$ws.Range("E28").Value = '  +6.79%  '
$ws.Range("D29").Value = '7.43'
$ws.Range("E29").Value = '  +10.93%  '
$ws.Range("D30").Value = '1.86'
$ws.Range("E30").Value = '  +6.58%  '
$ws.Range("D31").Value = '1.25'
$ws.Range("E31").Value = '  +6.78%  '
$ws.Range("D32").Value = '21.22'
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("D33").Value = '156.75'
$ws.Range("E33").Value = '  +2.16%  '
$ws.Range("E34").Value = '  +7.12%  '
$ws.Range("D35").Value = '6.02'
$ws.Range("E35").Value = '  +8.03%  '
$ws.Range("E36").Value = '  +4.08%  '
$ws.Range("D37").Value = '26.44'
$ws.Range("E37").Value = '  +14.32%  '
$ws.Range("D38").Value = '0.0705'
$ws.Range("E38").Value = '  +5.18%  '
$ws.Range("D39").Value = '3.095.18'
$ws.Range("E39").Value = '  +2.56%  '
$ws.Range("D40").Value = '37.92'
$ws.Range("E40").Value = '  +3.65%  '
$ws.Range("D41").Value = '3.92'
$ws.Range("E41").Value = '  +10.53%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").Value = '1.48'
$ws.Range("E43").Value = '  +5.73%  '
$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.665'
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("D45").Value = '2.337.28'
$ws.Range("E45").Value = '  +5.28%  '
$ws.Range("E46").Value = '  +4.47%  '
$ws.Range("D47").Value = '2.01'
$ws.Range("E47").Value = '  +4.06%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '6.05'
$ws.Range("E48").Value = '  +5.54%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0246'
$ws.Range("E49").Value = '  +3.65%  '
$ws.Range("D50").Value = '20.19'
$ws.Range("E50").Value = '  +6.29%  '
$ws.Range("D51").Value = '0.0903'
$ws.Range("E51").Value = '  +6.84%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
